$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @(2, 316828, 403794165),
    @(8, 853, 1254908),
    @(10, 116167, 170222627),
    @(12, 58691, 84702859),
    @(16, 3980, 5647432),
    @(19, 70, 101535),
    @(20, 6521, 9095322),
    @(22, 76552, 95507166),
    @(28, 32239, 47198204),
    @(30, 11373, 16358102),
    @(33, 1556, 2185307),
    @(35, 1788, 2523615),
    @(36, 96168, 121098094),
    @(42, 900, 1324685),
    @(44, 44089, 64613631),
    @(46, 9045, 12981599),
    @(48, 1397, 1940109),
    @(51, 2262, 3156394),
    @(52, 68303, 85704445),
    @(58, 27936, 40970319),
    @(61, 10972, 15864956),
    @(67, 1444, 2022065),
    @(69, 20279, 26561072),
    @(73, 7525, 11017030),
    @(75, 5060, 7347246),
    @(76, 485, 686239),
    @(77, 270, 378673),
    @(78, 139155, 173550046),
    @(80, 87, 121884),
    @(84, 63057, 92420284),
    @(87, 29398, 42523433),
    @(89, 2720, 3917242),
    @(90, 2771, 3914295),
    @(91, 32324, 43785476),
    @(95, 7825, 11506381),
    @(97, 7149, 10364336),
    @(100, 482, 695443),
    @(101, 8877, 12321466),
    @(103, 2235, 3292802),
    @(105, 3000, 4382451),
    @(108, 177, 251543),
    @(109, 139616, 172652934),
    @(115, 52270, 76626265),
    @(117, 26658, 38619596),
    @(121, 2210, 3103855),
    @(123, 495329, 653188045),
    @(128, 1366, 2025311),
    @(130, 205257, 301727346),
    @(133, 177451, 257923860),
    @(136, 2831, 3979381),
    @(138, 6201, 8761607),
    @(141, 43948, 58667555),
    @(147, 13929, 20427111),
    @(148, 3707, 5345180),
    @(151, 396, 569931),
    @(153, 376, 530751),
    @(154, 17317, 22880660),
    @(158, 7076, 10291203),
    @(160, 4925, 7086563),
    @(163, 263, 375864),
    @(165, 15613, 22651655),
    @(166, 1759, 2616230),
    @(167, 235, 347302),
    @(170, 80, 119949),
    @(171, 86726, 108486436),
    @(178, 33583, 49250791),
    @(180, 12858, 18575909),
    @(184, 1615, 2268162),
    @(186, 235809, 293150201),
    @(194, 85969, 126019371),
    @(197, 32693, 47050170),
    @(200, 5078, 7233437),
    @(203, 4777, 6612594),
    @(206, 260749, 322727090),
    @(207, 158, 173913),
    @(208, 252, 360587),
    @(215, 94410, 138116145),
    @(218, 50873, 73522453),
    @(221, 4649, 6527106),
    @(224, 5628, 7780678),
    @(227, 104992, 131367376),
    @(234, 49113, 71951585),
    @(236, 12241, 17599208),
    @(238, 1885, 2701382),
    @(240, 2447, 3419815),
    @(241, 254264, 321026342),
    @(249, 94952, 139133155),
    @(252, 64164, 92980807),
    @(254, 2394, 3377928),
    @(257, 4513, 6336039)
)

foreach ($u in $updates) {
    $row = $u[0]
    $nombreAides = $u[1]
    $montantTotal = $u[2]
    $ws.Cells.Item($row, 3).Value = $nombreAides
    $ws.Cells.Item($row, 4).Value = $montantTotal
}
